# Natmi following Dr Hou advice
# Updates the Fgf2-Sdc2 LR-pair sheet: recomputed statistics for the
# sending/target cluster combinations (ECs, FAPs, M1, sCs x ECs, FAPs, sCs),
# adding the new "sCs" sending-cluster rows (11-13) that were missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{row=2; A="ECs"; B="Fgf2"; C="Sdc2"; D="ECs"; E=2; F=0.6666666666666666; G=0.8896923333333334; H=2.669077; I=0.09326752961394506; J=0.09326752961394505; K=3; L=1; M=2.248835333333334; N=6.746506; O=0.03590294220158827; P=0.03590294220158827; Q=2.000771554995778; R=18.006943994962; S=0.003348578725014392; T=0.003348578725014392}
    @{row=3; A="ECs"; B="Fgf2"; C="Sdc2"; D="FAPs"; E=2; F=0.6666666666666666; G=0.8896923333333334; H=2.669077; I=0.09326752961394506; J=0.09326752961394505; K=3; L=1; M=44.29005966666667; N=132.870179; O=0.7070964373190639; P=0.7070964373190639; Q=39.40452652830923; R=354.6407387547831; S=0.06594913790757084; T=0.06594913790757083}
    @{row=4; A="ECs"; B="Fgf2"; C="Sdc2"; D="sCs"; E=2; F=0.6666666666666666; G=0.8896923333333334; H=2.669077; I=0.09326752961394506; J=0.09326752961394505; K=3; L=1; M=16.09762433333333; N=48.292873; O=0.2570006204793478; P=0.2570006204793479; Q=14.32193295424678; R=128.897396588221; S=0.02396981298135983; T=0.02396981298135983}
    @{row=5; A="FAPs"; B="Fgf2"; C="Sdc2"; D="ECs"; E=3; F=1; G=7.418580000000001; H=22.25574; I=0.777698766101638; J=0.777698766101638; K=3; L=1; M=2.248835333333334; N=6.746506; O=0.03590294220158827; P=0.03590294220158827; Q=16.68316482716001; R=150.14848344444; S=0.02792167384959363; T=0.02792167384959363}
    @{row=6; A="FAPs"; B="Fgf2"; C="Sdc2"; D="FAPs"; E=3; F=1; G=7.418580000000001; H=22.25574; I=0.777698766101638; J=0.777698766101638; K=3; L=1; M=44.29005966666667; N=132.870179; O=0.7070964373190639; P=0.7070964373190639; Q=328.5693508419401; R=2957.124157577461; S=0.5499080268179002; T=0.5499080268179002}
    @{row=7; A="FAPs"; B="Fgf2"; C="Sdc2"; D="sCs"; E=3; F=1; G=7.418580000000001; H=22.25574; I=0.777698766101638; J=0.777698766101638; K=3; L=1; M=16.09762433333333; N=48.292873; O=0.2570006204793478; P=0.2570006204793479; Q=119.42151392678; R=1074.79362534102; S=0.1998690654341441; T=0.1998690654341442}
    @{row=8; A="M1"; B="Fgf2"; C="Sdc2"; D="ECs"; E=3; F=1; G=0.298413; H=0.895239; I=0.031282997809377; J=0.03128299780937701; K=3; L=1; M=2.248835333333334; N=6.746506; O=0.03590294220158827; P=0.03590294220158827; Q=0.6710816983260001; R=6.039735284934; S=0.001123151662242475; T=0.001123151662242475}
    @{row=9; A="M1"; B="Fgf2"; C="Sdc2"; D="FAPs"; E=3; F=1; G=0.298413; H=0.895239; I=0.031282997809377; J=0.03128299780937701; K=3; L=1; M=44.29005966666667; N=132.870179; O=0.7070964373190639; P=0.7070964373190639; Q=13.216729575309; R=118.950566177781; S=0.02212009629967056; T=0.02212009629967056}
    @{row=10; A="M1"; B="Fgf2"; C="Sdc2"; D="sCs"; E=3; F=1; G=0.298413; H=0.895239; I=0.031282997809377; J=0.03128299780937701; K=3; L=1; M=16.09762433333333; N=48.292873; O=0.2570006204793478; P=0.2570006204793479; Q=4.803740370182999; R=43.233663331647; S=0.008039749847463968; T=0.008039749847463971}
    @{row=11; A="sCs"; B="Fgf2"; C="Sdc2"; D="ECs"; E=3; F=1; G=0.9324580000000001; H=2.797374; I=0.09775070647503986; J=0.09775070647503986; K=3; L=1; M=2.248835333333334; N=6.746506; O=0.03590294220158827; P=0.03590294220158827; Q=2.096944497249334; R=18.872500475244; S=0.003509537964737777; T=0.003509537964737777}
    @{row=12; A="sCs"; B="Fgf2"; C="Sdc2"; D="FAPs"; E=3; F=1; G=0.9324580000000001; H=2.797374; I=0.09775070647503986; J=0.09775070647503986; K=3; L=1; M=44.29005966666667; N=132.870179; O=0.7070964373190639; P=0.7070964373190639; Q=41.29862045666068; R=371.6875841099461; S=0.06911917629392224; T=0.06911917629392224}
    @{row=13; A="sCs"; B="Fgf2"; C="Sdc2"; D="sCs"; E=3; F=1; G=0.9324580000000001; H=2.797374; I=0.09775070647503986; J=0.09775070647503986; K=3; L=1; M=16.09762433333333; N=48.292873; O=0.2570006204793478; P=0.2570006204793479; Q=15.01035859061133; R=135.093227315502; S=0.02512199221637984; T=0.02512199221637985}
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rowsData) {
    foreach ($col in $columns) {
        $ws.Range($col + $r.row).Value = $r[$col]
    }
}
